$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 183, shifting old rows 183-186 down to 192-195
$ws.Range("A183:A191").EntireRow.Insert()

$appraiseB = 'Upon taking the Appraise skill, an artificer can select ONE expertise to be well-versed in (non-exhaustive examples in parentheses). Gems (Bangles, Broaches, Rings), Runes (Tattoos, Brands, Etchings), Weaponry (Blades, Hammers, Spears, etc.), Armor (Breastplates, Greaves, Helms, Shields, etc.), Spell Focus (Wizard Staves, Spellbooks, Holy Symbols, etc.), Alchemy (Potions, Reagents, Monster Parts, etc.), Machinery (Clockwork Machines, Engines, Traps, etc.), Clothing (Boots, Capes, Cloaks, Hats, etc.),Art (Paintings, Sculptures, Tapestries, etc.). The Artificer also gains the ability to write Schematics to replicate the effects of spells– embedding in Artificer Objects of your Appraisal expertise.  The Artificer may read scrolls, magical runes, and other materials as though they had the Read/Write Arcana spell, but only to understand them as a reference for creating items (see The Object’s Schematics below), not to cast the spells.  When scrolls are read only as a reference, they are neither activated nor expended. When using them to create a Schematic, the scroll or other reference material is expended due to experimentation.  The Schematic persists indefinitely for future projects, but must be at-hand to create an Artificer Object (they can be stolen through Pickpocket or Loot and are prioritized over coin).'
$appraiseE = 'In terms of roleplay, your keen eye for objects of this category allows you to intuit its relative worth and even altered properties, subject to Organizer oversight.  Make it fun and specific to your character!  You may take this skill a number of times equal to your Artificer Tier, each time selecting a new Expertise.  An Artificer can only understand reference materials in Tier up to their Artificer Professional Expertise level. An Artificer may read as many reference materials as they like, but may only draft one Schematic per event.  These Schematics are mundane and may not be used to cast spells like scrolls.  Once drafted, a Schematic does not scale—even if an Artificer that authored it improves, the Schematic stays the same level.'
$appraiseF = 'To Appraise - a magnifying glass, examiner’s loupe, or craftsperson-specific object.  For Schematics - an artistic architectural document that shows the fashioning of an object or demonstration of an ability.'

$ws.Range("A183").Value = 'Appraise [Gems]'
$ws.Range("B183").Value = $appraiseB
$ws.Range("C183").Value = 'Artificer'
$ws.Range("D183").Value = 1
$ws.Range("E183").Value = $appraiseE
$ws.Range("F183").Value = $appraiseF
$ws.Range("I183").Value = $false

$ws.Range("A184").Value = 'Appraise [Runes]'
$ws.Range("B184").Value = $appraiseB
$ws.Range("C184").Value = 'Artificer'
$ws.Range("D184").Value = 1
$ws.Range("E184").Value = $appraiseE
$ws.Range("F184").Value = $appraiseF
$ws.Range("I184").Value = $false

$ws.Range("A185").Value = 'Appraise [Weaponry]'
$ws.Range("B185").Value = $appraiseB
$ws.Range("C185").Value = 'Artificer'
$ws.Range("D185").Value = 1
$ws.Range("E185").Value = $appraiseE
$ws.Range("F185").Value = $appraiseF
$ws.Range("I185").Value = $false

$ws.Range("A186").Value = 'Appraise [Armor]'
$ws.Range("B186").Value = $appraiseB
$ws.Range("C186").Value = 'Artificer'
$ws.Range("D186").Value = 1
$ws.Range("E186").Value = $appraiseE
$ws.Range("F186").Value = $appraiseF
$ws.Range("I186").Value = $false

$ws.Range("A187").Value = 'Appraise [Spell Focus]'
$ws.Range("B187").Value = $appraiseB
$ws.Range("C187").Value = 'Artificer'
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = $appraiseE
$ws.Range("F187").Value = $appraiseF
$ws.Range("I187").Value = $false

$ws.Range("A188").Value = 'Appraise [Alchemy]'
$ws.Range("B188").Value = $appraiseB
$ws.Range("C188").Value = 'Artificer'
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = $appraiseE
$ws.Range("F188").Value = $appraiseF
$ws.Range("I188").Value = $false

$ws.Range("A189").Value = 'Appraise [Machinery]'
$ws.Range("B189").Value = $appraiseB
$ws.Range("C189").Value = 'Artificer'
$ws.Range("D189").Value = 1
$ws.Range("E189").Value = $appraiseE
$ws.Range("F189").Value = $appraiseF
$ws.Range("I189").Value = $false

$ws.Range("A190").Value = 'Appraise [Clothing]'
$ws.Range("B190").Value = $appraiseB
$ws.Range("C190").Value = 'Artificer'
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = $appraiseE
$ws.Range("F190").Value = $appraiseF
$ws.Range("I190").Value = $false

$ws.Range("A191").Value = 'Appraise [Art]'
$ws.Range("B191").Value = $appraiseB
$ws.Range("C191").Value = 'Artificer'
$ws.Range("D191").Value = 1
$ws.Range("E191").Value = $appraiseE
$ws.Range("F191").Value = $appraiseF
$ws.Range("I191").Value = $false

